$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8
$ws.Range("E8").Value = "Chief Data Officer"
$ws.Range("I8").Value = "Dependent on AI/ML Implementation milestone completion"
$ws.Range("J8").Value = "Critical action for Artificial Intelligence and Machine Learning success"

# Row 9
$ws.Range("E9").Value = "Data Scientists"
$ws.Range("I9").Value = "Dependent on AI/ML Implementation milestone completion"
$ws.Range("J9").Value = "Critical action for Artificial Intelligence and Machine Learning success"

# Row 10
$ws.Range("E10").Value = "ML Engineers"
$ws.Range("I10").Value = "Dependent on AI/ML Implementation milestone completion"
$ws.Range("J10").Value = "Critical action for Artificial Intelligence and Machine Learning success"

# Row 11
$ws.Range("E11").Value = "Business Analysts"
$ws.Range("I11").Value = "Dependent on AI/ML Implementation milestone completion"
$ws.Range("J11").Value = "Critical action for Artificial Intelligence and Machine Learning success"

# Row 12 (Owner unchanged)
$ws.Range("I12").Value = "Dependent on AI/ML Implementation milestone completion"
$ws.Range("J12").Value = "Critical action for Artificial Intelligence and Machine Learning success"

# Row 13 (Owner unchanged)
$ws.Range("I13").Value = "Dependent on AI/ML Implementation milestone completion"
$ws.Range("J13").Value = "Critical action for Artificial Intelligence and Machine Learning success"

# Row 14
$ws.Range("E14").Value = "Chief Data Officer"
$ws.Range("I14").Value = "Dependent on AI/ML Implementation milestone completion"
$ws.Range("J14").Value = "Critical action for Artificial Intelligence and Machine Learning success"

# Row 15
$ws.Range("E15").Value = "Data Scientists"
$ws.Range("I15").Value = "Dependent on AI/ML Implementation milestone completion"
$ws.Range("J15").Value = "Critical action for Artificial Intelligence and Machine Learning success"

# Row 16
$ws.Range("E16").Value = "ML Engineers"
$ws.Range("I16").Value = "Dependent on AI/ML Implementation milestone completion"
$ws.Range("J16").Value = "Critical action for Artificial Intelligence and Machine Learning success"

# Row 17
$ws.Range("E17").Value = "Business Analysts"
$ws.Range("I17").Value = "Dependent on AI/ML Implementation milestone completion"
$ws.Range("J17").Value = "Critical action for Artificial Intelligence and Machine Learning success"
